$d = $word.ActiveDocument
$ndash = [char]0x2013

# ---------------------------------------------------------------------------
# Paragraph 1 ("Research question: ...") — originally a single run holding
# the whole sentence, so it is replaced in one shot.
# ---------------------------------------------------------------------------
$old1 = "Research question: Is there any difference in mean of Temperature between Year 1961" + $ndash + "1981, Year 1982" + $ndash + "2001 and Year 2002" + $ndash + "2022? "
$new1 = "Research question: Is there any difference in median of Global  Average Temperature anamoly between the Periods 1961 " + $ndash + " 1991 and 1992 - 2022?"

$r1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Output "Para1 replace: $r1"

# ---------------------------------------------------------------------------
# Paragraph 2 ("RQ: Is there any difference in mean of Global  Average
# Temperature between Year 1961 -1981, ..."). This paragraph already carries
# a "Global  Average" run wrapped in proofErr gramStart/gramEnd markers that
# is untouched by the edit, so the replacement is split into two calls that
# do not cross into that run (or the unchanged " Temperature " run after
# it), to keep that markup intact.
# ---------------------------------------------------------------------------

# 1) "in mean of " -> "in median of " (fully before "Global  Average",
#    scoped to the "RQ:" paragraph only via a narrow, unique anchor so the
#    untouched "Global  Average" run/proofErr markers right after it are
#    left completely alone)
$r2a = $d.Content.Find.Execute("difference in mean of ", $true, $false, $false, $false, $false, $true, 1, $false, "difference in median of ", 2)
Write-Output "Para2 mean->median replace: $r2a"

# 2) "between Year 1961 -1981, Year 1982 - 2001 and Year 2002 - 2022? "
#    -> "anamoly between the Periods 1961 - 1991 and 1992 - 2022?"
#    (fully after the unchanged "Global  Average" / " Temperature " runs)
$old2b = "between Year 1961 -1981, Year 1982 " + $ndash + " 2001 and Year 2002 - 2022? "
$new2b = "anamoly between the Periods 1961 " + $ndash + " 1991 and 1992 - 2022?"

$r2b = $d.Content.Find.Execute($old2b, $true, $false, $false, $false, $false, $true, 1, $false, $new2b, 2)
Write-Output "Para2 between-dates replace: $r2b"
